$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 708.1177
$ws.Range("J17").Value = 708.1177
$ws.Range("L17").Value = 2124.3531
$ws.Range("N17").Value = -2460.3531
$ws.Range("H29").Value = 3178.5715
$ws.Range("J29").Value = 4087.5
$ws.Range("L29").Value = 12262.5
$ws.Range("N29").Value = -12824.5
$ws.Range("H32").Value = 9027.286
$ws.Range("I32").Value = 9158.200000000001
$ws.Range("K32").Value = 9158.200000000001
$ws.Range("M32").Value = -8832.200000000001
$ws.Range("H38").Value = 200
$ws.Range("I38").Value = 200
$ws.Range("K38").Value = 600
$ws.Range("M38").Value = -228
$ws.Range("H40").Value = 766.8333
$ws.Range("I40").Value = 775.25
$ws.Range("K40").Value = 775.25
$ws.Range("M40").Value = -600.25
$ws.Range("I62").Value = 6249.5
$ws.Range("J62").Value = 5900
$ws.Range("K62").Value = 6249.5
$ws.Range("L62").Value = 5900
$ws.Range("M62").Value = -5625.5
$ws.Range("N62").Value = -7148
$ws.Range("I65").Value = 6249.5
$ws.Range("J65").Value = 5900
$ws.Range("K65").Value = 31247.5
$ws.Range("L65").Value = 29500
$ws.Range("M65").Value = -28127.5
$ws.Range("N65").Value = -35740
$ws.Range("H98").Value = 826.25
$ws.Range("I98").Value = 826.25
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 826.25
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 671.75
$ws.Range("N98").ClearContents()
$ws.Range("H107").Value = 1222.5714
$ws.Range("I107").Value = 1201.2307
$ws.Range("K107").Value = 1201.2307
$ws.Range("M107").Value = 718.7692999999999
$ws.Range("H122").Value = 826.25
$ws.Range("I122").Value = 826.25
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2478.75
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -28.75
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 1302.9131
$ws.Range("I132").Value = 1320.9546
$ws.Range("K132").Value = 3962.8638
$ws.Range("M132").Value = -1432.8638
$ws.Range("H137").Value = 3319.0667
$ws.Range("I137").Value = 3048.3333
$ws.Range("J137").Value = 3499.5557
$ws.Range("K137").Value = 9144.999899999999
$ws.Range("L137").Value = 10498.6671
$ws.Range("M137").Value = -6594.999899999999
$ws.Range("N137").Value = -15598.6671

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 465.75
$ws.Range("I2").Value = 465.75
$ws.Range("K2").Value = 465.75
$ws.Range("M2").Value = -352.75
$ws.Range("H32").Value = 3916.04
$ws.Range("I32").Value = 3430.9565
$ws.Range("K32").Value = 3430.9565
$ws.Range("M32").Value = -3143.9565
$ws.Range("H45").Value = 2507.5
$ws.Range("I45").Value = 2394.1667
$ws.Range("K45").Value = 2394.1667
$ws.Range("M45").Value = -2017.1667
$ws.Range("H61").Value = 3468.5
$ws.Range("J61").Value = 2518.8
$ws.Range("L61").Value = 2518.8
$ws.Range("N61").Value = -2942.8
$ws.Range("H102").Value = 1370.875
$ws.Range("I102").Value = 1503
$ws.Range("K102").Value = 1503
$ws.Range("M102").Value = 119
$ws.Range("H110").Value = 4606.3335
$ws.Range("I110").Value = 4606.3335
$ws.Range("K110").Value = 4606.3335
$ws.Range("M110").Value = -2561.3335
$ws.Range("H116").Value = 465.75
$ws.Range("I116").Value = 465.75
$ws.Range("K116").Value = 465.75
$ws.Range("M116").Value = 1828.25
$ws.Range("H122").Value = 1572.8
$ws.Range("I122").Value = 1572.8
$ws.Range("K122").Value = 4718.4
$ws.Range("M122").Value = -2268.4
$ws.Range("H132").Value = 2076.516
$ws.Range("I132").Value = 1542.5454
$ws.Range("K132").Value = 4627.6362
$ws.Range("M132").Value = -2097.6362
$ws.Range("H136").Value = 3468.5
$ws.Range("J136").Value = 2518.8
$ws.Range("L136").Value = 7556.400000000001
$ws.Range("N136").Value = -12656.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 465.75
$ws.Range("I3").Value = 465.75
$ws.Range("K3").Value = 465.75
$ws.Range("M3").Value = -351.75
$ws.Range("H107").Value = 4713
$ws.Range("I107").Value = 4709.2104
$ws.Range("K107").Value = 4709.2104
$ws.Range("M107").Value = -2789.2104

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 200000
$ws.Range("J9").Value = 200000
$ws.Range("L9").Value = 200000
$ws.Range("N9").Value = -200336
$ws.Range("H31").Value = 4938.074
$ws.Range("I31").Value = 1598.6154
$ws.Range("J31").Value = 8039
$ws.Range("K31").Value = 1598.6154
$ws.Range("L31").Value = 8039
$ws.Range("M31").Value = -1303.6154
$ws.Range("N31").Value = -8629
$ws.Range("H34").Value = 4938.074
$ws.Range("I34").Value = 1598.6154
$ws.Range("J34").Value = 8039
$ws.Range("K34").Value = 1598.6154
$ws.Range("L34").Value = 8039
$ws.Range("M34").Value = -1396.6154
$ws.Range("N34").Value = -8443
$ws.Range("H58").Value = 2578
$ws.Range("I58").Value = 2497.3333
$ws.Range("J58").Value = 2699
$ws.Range("K58").Value = 2497.3333
$ws.Range("L58").Value = 2699
$ws.Range("M58").Value = -2294.3333
$ws.Range("N58").Value = -3105
$ws.Range("H99").Value = 2499.6667
$ws.Range("I99").Value = 2499.6667
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 2499.6667
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -1001.6667
$ws.Range("N99").ClearContents()
$ws.Range("H126").Value = 2499.6667
$ws.Range("I126").Value = 2499.6667
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 7499.000100000001
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -5029.000100000001
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 2801.5789
$ws.Range("I132").Value = 1717.2222
$ws.Range("K132").Value = 5151.6666
$ws.Range("M132").Value = -2621.6666
$ws.Range("H134").Value = 3264.1667
$ws.Range("I134").Value = 3117.1
$ws.Range("J134").Value = 3999.5
$ws.Range("K134").Value = 9351.299999999999
$ws.Range("L134").Value = 11998.5
$ws.Range("M134").Value = -6816.299999999999
$ws.Range("N134").Value = -17068.5
$ws.Range("H136").Value = 2578
$ws.Range("I136").Value = 2497.3333
$ws.Range("J136").Value = 2699
$ws.Range("K136").Value = 7491.999899999999
$ws.Range("L136").Value = 8097
$ws.Range("M136").Value = -4941.999899999999
$ws.Range("N136").Value = -13197

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1282.375
$ws.Range("J68").Value = 1294.1428
$ws.Range("L68").Value = 3882.4284
$ws.Range("N68").Value = -5504.428400000001
$ws.Range("H71").Value = 1282.375
$ws.Range("J71").Value = 1294.1428
$ws.Range("L71").Value = 11647.2852
$ws.Range("N71").Value = -19759.2852
$ws.Range("H113").Value = 839.0909
$ws.Range("I113").Value = 749
$ws.Range("J113").Value = 859.1111
$ws.Range("K113").Value = 2247
$ws.Range("L113").Value = 2577.3333
$ws.Range("M113").Value = -77
$ws.Range("N113").Value = -6917.3333
$ws.Range("H129").Value = 744.5
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 5058.3335
$ws.Range("J3").Value = 6498.75
$ws.Range("L3").Value = 6498.75
$ws.Range("N3").Value = -6730.75
$ws.Range("H122").Value = 1633.1666
$ws.Range("I122").Value = 1633.1666
$ws.Range("K122").Value = 4899.4998
$ws.Range("M122").Value = -2449.4998
$ws.Range("H132").Value = 3184.147
$ws.Range("I132").Value = 2702.1155
$ws.Range("K132").Value = 8106.3465
$ws.Range("M132").Value = -5576.3465

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3184.8125
$ws.Range("I40").Value = 3121.3076
$ws.Range("J40").Value = 3460
$ws.Range("K40").Value = 3121.3076
$ws.Range("L40").Value = 3460
$ws.Range("M40").Value = -2985.3076
$ws.Range("N40").Value = -3732
$ws.Range("H55").Value = 585
$ws.Range("I55").Value = 540
$ws.Range("K55").Value = 540
$ws.Range("M55").Value = -367
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H93").Value = 668
$ws.Range("I93").Value = 750
$ws.Range("J93").Value = 504
$ws.Range("K93").Value = 750
$ws.Range("L93").Value = 504
$ws.Range("M93").Value = 498
$ws.Range("N93").Value = -3000
$ws.Range("H133").Value = 140000
$ws.Range("J133").Value = 140000
$ws.Range("L133").Value = 140000
$ws.Range("N133").Value = -145060

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 30000
$ws.Range("J63").Value = 30000
$ws.Range("L63").Value = 30000
$ws.Range("N63").Value = -31248
$ws.Range("H66").Value = 30000
$ws.Range("J66").Value = 30000
$ws.Range("L66").Value = 90000
$ws.Range("N66").Value = -96240
$ws.Range("H81").Value = 8152.1816
$ws.Range("I81").Value = 5083.75
$ws.Range("K81").Value = 10167.5
$ws.Range("M81").Value = -9106.5
$ws.Range("H84").Value = 8152.1816
$ws.Range("I84").Value = 5083.75
$ws.Range("K84").Value = 50837.5
$ws.Range("M84").Value = -45533.5
$ws.Range("H126").Value = 2030.375
$ws.Range("I126").Value = 1963.2858
$ws.Range("K126").Value = 5889.857400000001
$ws.Range("M126").Value = -3419.857400000001
$ws.Range("H136").Value = 8768
$ws.Range("I136").Value = 12878.6
$ws.Range("K136").Value = 38635.8
$ws.Range("M136").Value = -36085.8
